$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.319.75"
$ws.Range("E2").Value = "  +2.36%  "
$ws.Range("D3").Value = "1.820.29"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4651"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3769"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07431"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8717"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "1.823.23"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.684"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.409"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07082"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008751"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").Value = "27.329.79"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.310"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("D24").Value = "2.049.94"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.252"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.294"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08896"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7813"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.183"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.527"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.928"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.099"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01971"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05255"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.289"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5298"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.367"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +20.54%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.895"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.612"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5043"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.675"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06327"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.64%  "
